{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"5840560 - Marco Antonio Carvalho Pereira\";\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text.trim() === targetText) {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$target = \"5840560 - Marco Antonio Carvalho Pereira\"\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Trim() -eq $target) {\n        $p.Range.Delete()\n    }\n}\n"}
